$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 55
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
